# Scheduled-runner refresh: updates market-price-derived leve profit
# columns (currentAveragePrice*, LevePrice*, LeveProfit*) on several
# rows across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 8560.286
$ws.Range("J62").Value = 11372.5
$ws.Range("L62").Value = 11372.5
$ws.Range("N62").Value = -12620.5

$ws.Range("H65").Value = 8560.286
$ws.Range("J65").Value = 11372.5
$ws.Range("L65").Value = 56862.5
$ws.Range("N65").Value = -63102.5

$ws.Range("H80").Value = 485.4516
$ws.Range("I80").Value = 318.15384
$ws.Range("J80").Value = 606.2778
$ws.Range("K80").Value = 954.4615200000001
$ws.Range("L80").Value = 1818.8334
$ws.Range("M80").Value = 43.53847999999994
$ws.Range("N80").Value = -3814.8334

$ws.Range("H83").Value = 485.4516
$ws.Range("I83").Value = 318.15384
$ws.Range("J83").Value = 606.2778
$ws.Range("K83").Value = 2863.38456
$ws.Range("L83").Value = 5456.500199999999
$ws.Range("M83").Value = 2128.61544
$ws.Range("N83").Value = -15440.5002

$ws.Range("H100").Value = 5748.5
$ws.Range("I100").Value = 4001
$ws.Range("K100").Value = 4001
$ws.Range("M100").Value = -3460

$ws.Range("H111").Value = 1207.5714
$ws.Range("I111").Value = 1238.25
$ws.Range("J111").Value = 1166.6666
$ws.Range("K111").Value = 3714.75
$ws.Range("L111").Value = 3499.9998
$ws.Range("M111").Value = -647.75
$ws.Range("N111").Value = -9633.9998

$ws.Range("H112").Value = 1665
$ws.Range("I112").Value = 1665
$ws.Range("K112").Value = 4995
$ws.Range("M112").Value = -3887

$ws.Range("H113").Value = 11756.143
$ws.Range("J113").Value = 13947
$ws.Range("L113").Value = 13947
$ws.Range("N113").Value = -20455

$ws.Range("H125").Value = 2910
$ws.Range("I125").Value = 2620.9092
$ws.Range("K125").Value = 23588.1828
$ws.Range("M125").Value = -21128.1828

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3455548
$ws.Range("I32").Value = 4277.826
$ws.Range("J32").Value = 16685416
$ws.Range("K32").Value = 4277.826
$ws.Range("L32").Value = 16685416
$ws.Range("M32").Value = -3990.826
$ws.Range("N32").Value = -16685990

$ws.Range("H102").Value = 9201.799999999999
$ws.Range("I102").Value = 6009
$ws.Range("K102").Value = 6009
$ws.Range("M102").Value = -4387

$ws.Range("H122").Value = 2715.8572
$ws.Range("I122").Value = 2715.8572
$ws.Range("K122").Value = 8147.571599999999
$ws.Range("M122").Value = -5697.571599999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H13").Value = 49999
$ws.Range("I13").Value = 49999
$ws.Range("J13").Value = 0
$ws.Range("K13").Value = 49999
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = -49831
$ws.Range("N13").Value = $null

$ws.Range("H22").Value = 579.4
$ws.Range("I22").Value = 498.5
$ws.Range("K22").Value = 498.5
$ws.Range("M22").Value = -325.5

$ws.Range("H63").Value = 100000
$ws.Range("J63").Value = 100000
$ws.Range("L63").Value = 100000
$ws.Range("N63").Value = -101372

$ws.Range("H66").Value = 100000
$ws.Range("J66").Value = 100000
$ws.Range("L66").Value = 300000
$ws.Range("N66").Value = -306864

$ws.Range("H86").Value = 6318.364
$ws.Range("I86").Value = 3336.6667
$ws.Range("J86").Value = 7436.5
$ws.Range("K86").Value = 3336.6667
$ws.Range("L86").Value = 7436.5
$ws.Range("M86").Value = -2213.6667
$ws.Range("N86").Value = -9682.5

$ws.Range("H89").Value = 6318.364
$ws.Range("I89").Value = 3336.6667
$ws.Range("J89").Value = 7436.5
$ws.Range("K89").Value = 16683.3335
$ws.Range("L89").Value = 37182.5
$ws.Range("M89").Value = -11067.3335
$ws.Range("N89").Value = -48414.5

$ws.Range("H99").Value = 3635.3
$ws.Range("I99").Value = 3874.8333
$ws.Range("J99").Value = 3276
$ws.Range("K99").Value = 3874.8333
$ws.Range("L99").Value = 3276
$ws.Range("M99").Value = -2376.8333
$ws.Range("N99").Value = -6272

$ws.Range("H105").Value = 1770.3077
$ws.Range("I105").Value = 1626.875
$ws.Range("J105").Value = 1999.8
$ws.Range("K105").Value = 1626.875
$ws.Range("L105").Value = 1999.8
$ws.Range("M105").Value = 120.125
$ws.Range("N105").Value = -5493.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 5000
$ws.Range("I86").Value = 5000
$ws.Range("K86").Value = 5000
$ws.Range("M86").Value = -3877

$ws.Range("H89").Value = 5000
$ws.Range("I89").Value = 5000
$ws.Range("K89").Value = 25000
$ws.Range("M89").Value = -19384

$ws.Range("H116").Value = 76000
$ws.Range("J116").Value = 76000
$ws.Range("L116").Value = 76000
$ws.Range("N116").Value = -85178

$ws.Range("H141").Value = 85084
$ws.Range("J141").Value = 85084
$ws.Range("L141").Value = 85084
$ws.Range("N141").Value = -95444

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 1666.6666

$ws.Range("H131").Value = 3264.1538
$ws.Range("I131").Value = 2016.875
$ws.Range("K131").Value = 6050.625
$ws.Range("M131").Value = -1010.625

$ws.Range("H137").Value = 1499
$ws.Range("I137").Value = 1499
$ws.Range("K137").Value = 4497
$ws.Range("M137").Value = 603

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 4454.8184
$ws.Range("I55").Value = 4988.1665
$ws.Range("K55").Value = 4988.1665
$ws.Range("M55").Value = -4661.1665

$ws.Range("H132").Value = 47946
$ws.Range("I132").Value = 72318.47
$ws.Range("J132").Value = 7325.222
$ws.Range("K132").Value = 216955.41
$ws.Range("L132").Value = 21975.666
$ws.Range("M132").Value = -214425.41
$ws.Range("N132").Value = -27035.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 38011.5
$ws.Range("I45").Value = 16000
$ws.Range("K45").Value = 16000
$ws.Range("M45").Value = -15593

$ws.Range("H46").Value = 7137.357
$ws.Range("I46").Value = 6385.6
$ws.Range("K46").Value = 6385.6
$ws.Range("M46").Value = -6197.6

$ws.Range("H56").Value = 15703.143
$ws.Range("I56").Value = 11984.6
$ws.Range("K56").Value = 11984.6
$ws.Range("M56").Value = -11293.6

$ws.Range("H61").Value = 4826.0835
$ws.Range("I61").Value = 2566.1428
$ws.Range("K61").Value = 2566.1428
$ws.Range("M61").Value = -2364.1428

$ws.Range("H62").Value = 20239
$ws.Range("J62").Value = 20239
$ws.Range("L62").Value = 20239
$ws.Range("N62").Value = -21487

$ws.Range("H65").Value = 20239
$ws.Range("J65").Value = 20239
$ws.Range("L65").Value = 60717
$ws.Range("N65").Value = -66957

$ws.Range("H113").Value = 4826.0835
$ws.Range("I113").Value = 2566.1428
$ws.Range("K113").Value = 2566.1428
$ws.Range("M113").Value = -396.1428000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 27361.334
$ws.Range("I58").Value = 27361.334
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 27361.334
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -27053.334
$ws.Range("N58").Value = $null

$ws.Range("H62").Value = 12000
$ws.Range("J62").Value = 12000
$ws.Range("L62").Value = 12000
$ws.Range("N62").Value = -13248

$ws.Range("H65").Value = 12000
$ws.Range("J65").Value = 12000
$ws.Range("L65").Value = 60000
